$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117, shifting rows 117:233 down to 118:234
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new record's data
$ws.Cells.Item(117, 1).Value = 10
$ws.Cells.Item(117, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(117, 3).Value = "La Araucanía"
$ws.Cells.Item(117, 4).Value = 44601
$ws.Cells.Item(117, 5).Value = 9
$ws.Cells.Item(117, 6).Value = 100112001
$ws.Cells.Item(117, 7).Value = "Berenjena"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 25
$ws.Cells.Item(117, 11).Value = 12000
$ws.Cells.Item(117, 12).Value = 12000
$ws.Cells.Item(117, 13).Value = 12000
$ws.Cells.Item(117, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(117, 15).Value = "Región del Maule"
$ws.Cells.Item(117, 16).Value = 200
$ws.Cells.Item(117, 17).Value = 60
$ws.Cells.Item(117, 18).Value = "Hortaliza"
